# Weekly update: add two new daily-price records for Cilantro at
# "Terminal Hortofrutícola Agro Chillán", inserted as the first two data
# rows (row 38 and 39), pushing all existing records down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right above the current first data row (row 38).
# Doing this twice at the same index pushes everything below it down by 2.
$ws.Rows.Item(38).Insert()
$ws.Rows.Item(38).Insert()

# Common columns shared by every record in this sheet.
$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$codreg    = 16
$catId     = 100112040
$categoria = "Cilantro"
$variedad  = "Sin especificar"
$unidad    = "$/atado 0,5 a 1 kilo"
$origen    = "Provincia de Diguillín"
$kgUnid    = 1
$clasif    = "Hortaliza"

# New row 38: Cilantro "Primera", fecha 2022-08-17 (serial 44790)
$r = 38
$ws.Cells.Item($r,1).Value  = $mercadoId
$ws.Cells.Item($r,2).Value  = $mercado
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 44790
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $catId
$ws.Cells.Item($r,7).Value  = $categoria
$ws.Cells.Item($r,8).Value  = $variedad
$ws.Cells.Item($r,9).Value  = "Primera"
$ws.Cells.Item($r,10).Value = 200
$ws.Cells.Item($r,11).Value = 700
$ws.Cells.Item($r,12).Value = 800
$ws.Cells.Item($r,13).Value = 750
$ws.Cells.Item($r,14).Value = $unidad
$ws.Cells.Item($r,15).Value = $origen
$ws.Cells.Item($r,16).Value = 750
$ws.Cells.Item($r,17).Value = $kgUnid
$ws.Cells.Item($r,18).Value = $clasif

# New row 39: Cilantro "Segunda", fecha 2022-08-17 (serial 44790)
$r = 39
$ws.Cells.Item($r,1).Value  = $mercadoId
$ws.Cells.Item($r,2).Value  = $mercado
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 44790
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $catId
$ws.Cells.Item($r,7).Value  = $categoria
$ws.Cells.Item($r,8).Value  = $variedad
$ws.Cells.Item($r,9).Value  = "Segunda"
$ws.Cells.Item($r,10).Value = 160
$ws.Cells.Item($r,11).Value = 600
$ws.Cells.Item($r,12).Value = 600
$ws.Cells.Item($r,13).Value = 600
$ws.Cells.Item($r,14).Value = $unidad
$ws.Cells.Item($r,15).Value = $origen
$ws.Cells.Item($r,16).Value = 600
$ws.Cells.Item($r,17).Value = $kgUnid
$ws.Cells.Item($r,18).Value = $clasif

# Keep the date column's number format consistent with the rest of column D.
$ws.Range("D38:D39").NumberFormat = $ws.Range("D40").NumberFormat
